$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Intro GRP paragraph - reword
# ------------------------------------------------------------------
$oldIntro = "The Gross Regional Product (GRP) of San Diego County has demonstrated a consistent upward trajectory from 2019 to 2023, reflecting a robust economic environment. In 2019, the total GRP was approximately `$244.28 billion. Despite the challenges posed by the COVID-19 pandemic, 2020 saw a slight increase to `$244.82 billion. This resilience set the stage for a significant economic recovery in 2021, with the GRP rising to `$268.87 billion. The growth momentum continued into 2022, reaching `$296.68 billion, and further increased to `$308.71 billion in 2023, indicating sustained economic expansion."
$newIntro = "The Gross Regional Product (GRP) serves as a vital economic indicator, reflecting the economic output of San Diego County. Over the period from 2019 to 2023, the county has demonstrated remarkable economic resilience and growth. In 2019, the total GRP was approximately `$244.28 billion. Despite the challenges posed by the COVID-19 pandemic, 2020 saw a slight increase to `$244.82 billion. This resilience set the stage for a significant economic recovery in 2021, with the GRP rising to `$268.87 billion. The upward trend continued in 2022, reaching `$296.68 billion, and further increased to `$308.71 billion in 2023, indicating sustained economic growth."
$d.Content.Find.Execute($oldIntro, $true, $false, $false, $false, $false, $true, 1, $false, $newIntro, 2) | Out-Null

# ------------------------------------------------------------------
# 2) Per-capita paragraph becomes the chart-error placeholder
#    (set directly on the paragraph Range rather than via Find/Replace
#    so that the straight apostrophes in "'scales'" are not mangled
#    into curly quotes by autoformatting).
# ------------------------------------------------------------------
$newPerCapitaPlaceholder = "ERROR GENERATING CHART: 'scales'"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Per capita GRP also followed a positive trend")) {
        $p.Range.Text = $newPerCapitaPlaceholder
        break
    }
}

# ------------------------------------------------------------------
# 3) Industry paragraph + picture + caption collapse into the new
#    per-capita paragraph.
# ------------------------------------------------------------------
$newPerCapitaText = "The per capita GRP also reflects this positive economic trajectory. In 2019, the per capita GRP was approximately `$73,347. It increased to `$74,278 in 2020, despite the pandemic's impact. The recovery in 2021 was marked by a significant rise to `$82,100. This growth continued in 2022, with the per capita GRP reaching `$90,557, and further increased to `$94,916 in 2023. These figures underscore the robust economic health and increasing productivity of the region."

# Find the paragraph that currently holds the industry-contribution text
# (it precedes the inline picture + "Industry Contributions..." caption)
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("San Diego County's diverse economy")) {
        $targetPara = $i
        break
    }
}

if ($targetPara -ne $null) {
    # Replace its text with the reworded per-capita paragraph ...
    $d.Paragraphs($targetPara).Range.Text = $newPerCapitaText
    # ... then remove the following picture paragraph and the
    # "Industry Contributions to GRP (2023)" caption paragraph that
    # used to sit right after it.
    $d.Paragraphs($targetPara + 1).Range.Delete() | Out-Null
    $d.Paragraphs($targetPara + 1).Range.Delete() | Out-Null
}

# ------------------------------------------------------------------
# 4) Table: Industry/Contribution (11 rows x 2 cols) becomes
#    Year/Total GRP/Per Capita GRP (6 rows x 3 cols)
# ------------------------------------------------------------------
$t = $d.Tables(1)
$t.Columns.Add() | Out-Null

# Remove the five industry rows beyond "Finance and Insurance" (rows 7-11)
for ($i = 0; $i -lt 5; $i++) {
    $t.Rows(7).Delete() | Out-Null
}

$tableData = @(
    @("Year", "Total GRP (Billion `$)", "Per Capita GRP (`$)"),
    @("2019", "244.28", "73,347"),
    @("2020", "244.82", "74,278"),
    @("2021", "268.87", "82,100"),
    @("2022", "296.68", "90,557"),
    @("2023", "308.71", "94,916")
)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $t.Cell($r, $c).Range.Text = $tableData[$r - 1][$c - 1]
    }
}

for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $t.Columns($c).Width = 144
}

# ------------------------------------------------------------------
# 5) Caption after the table gets new text, plus a brand-new trailing
#    paragraph describing the 2023 industry breakdown.
# ------------------------------------------------------------------
$oldCaption = "Detailed Industry Contributions to GRP (2023)"
$newCaption = "San Diego County GRP and Per Capita GRP (2019-2023)"
$d.Content.Find.Execute($oldCaption, $true, $false, $false, $false, $false, $true, 1, $false, $newCaption, 2) | Out-Null

$industryParagraph = "In 2023, the contributions of various industries to the GRP highlight the diverse economic landscape of San Diego County. The government sector emerged as the largest contributor, with approximately `$52.92 billion. This was followed by the professional, scientific, and technical services sector, contributing around `$37.04 billion. Manufacturing added approximately `$31.67 billion, while the health care and social assistance sector contributed `$20.21 billion. The finance and insurance industry added `$19.50 billion, and the information sector contributed `$14.90 billion. Real estate and rental and leasing, retail trade, accommodation and food services, and administrative and support and waste management and remediation services also made significant contributions, underscoring the county's economic diversity and strength."

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = $industryParagraph
$newPara.Style = "Normal"

# ------------------------------------------------------------------
# 6) Header/footer updates
# ------------------------------------------------------------------
$sec = $d.Sections(1)

$header = $sec.Headers(1)
$header.Range.Find.Execute("San Diego County Economic Report", $true, $false, $false, $false, $false, $true, 1, $false, "Economic Growth and Industry Contributions in San Diego County", 2) | Out-Null

$footer = $sec.Footers(1)
$footer.Range.Find.Execute("Page \# of \#", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
